# Handback report regeneration: refresh the "Latest HO Xliff Generate Date"
# on the Overview sheet and the "Correspond Handoff/Handback DateTime"
# columns on the per-locale sheets for the file that was just handed back
# (2e797867-f734-4fca-9c49-cc639c1bfddf.md).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 2 is the 2e797867 file, row 3 is the af136f81 file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 18:50:09"
$wsOverview.Range("G3").Value = "2016-08-25 18:49:06"

# --- zh-cn sheet: column H = Correspond Handoff Datetime, column K = Correspond Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 18:49:58"
$wsZhCn.Range("K2").Value = "2016-08-25 18:50:33"
$wsZhCn.Range("H3").Value = "2016-08-25 18:48:57"
$wsZhCn.Range("K3").Value = "2016-08-25 18:49:31"

# --- de-de sheet: column H = Correspond Handoff Datetime, column K = Correspond Handback DateTime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 18:50:09"
$wsDeDe.Range("K2").Value = "2016-08-25 18:50:40"
$wsDeDe.Range("H3").Value = "2016-08-25 18:49:06"
$wsDeDe.Range("K3").Value = "2016-08-25 18:49:38"
